# The deck shipped with two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" / "Office" colour scheme
#   ppt/theme/theme2.xml -> "Integral"     / "Red Violet" colour scheme
# (theme2.xml is the one actually wired to the slide master / presentation,
# so it is what drives the look of the deck; theme1.xml is only referenced
# by the notes master.)
#
# The authored edit swaps the two themes' contents: the deck's live theme
# becomes the "Office Theme" palette (and the notes-only theme becomes the
# "Integral" palette). Reproduce that by re-pointing the live colour scheme
# (SlideMaster.ColorScheme, the PowerPoint OM's handle onto the active
# theme's <a:clrScheme>) at the Office palette's twelve colours, in the
# fixed PowerPoint colour-scheme order: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme

function HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officePalette = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 0; $i -lt $officePalette.Length; $i++) {
    $cs.Colors($i + 1).RGB = HexToRgb $officePalette[$i]
}
